# Adds two new paragraphs ("Hjkhkijij" and "ggiu") after the existing
# "This will be part of git project" paragraph, and moves the hidden
# "_GoBack" bookmark (left behind by the last edit) from the end of that
# paragraph to the end of the new last paragraph ("ggiu").

$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits at the end of the last paragraph
# ("This will be part of git project"). Remove it from there - it will be
# re-added after the new content, mirroring where Word leaves it following
# the most recent edit.
$goBack = $d.Bookmarks("_GoBack")
[void]$goBack.Delete()

# New paragraph: "Hjkhkijij" (flagged by the spell checker, hence the
# surrounding proofErr markers).
$lastPara = $d.Paragraphs.Last
[void]$lastPara.Range.InsertParagraphAfter()
$hjkPara = $d.Paragraphs.Last
[void]$hjkPara.Range.InsertXML(
    "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>Hjkhkijij</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "</w:p>"
)

# New paragraph: "ggiu", with the "_GoBack" bookmark collapsed right after
# the text (before the paragraph mark) - mirroring where it originally sat
# after the previous paragraph's text. The bookmark is written directly as
# part of the inserted markup (rather than via Bookmarks.Add) because this
# host mis-resolves a zero-length Range sitting immediately before a
# paragraph mark.
$hjkPara = $d.Paragraphs.Last
[void]$hjkPara.Range.InsertParagraphAfter()
$ggiuPara = $d.Paragraphs.Last
[void]$ggiuPara.Range.InsertXML(
    "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
    "<w:r><w:t>ggiu</w:t></w:r>" +
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" +
    "<w:bookmarkEnd w:id='0'/>" +
    "</w:p>"
)
